# Update "想去人数" (number of people interested) figures in the F column
# for sheets "展览" and "全部类型" (sheet1 and sheet4), leaving the other
# sheets ("演出" and "本地生活", which only contain header rows) untouched.

$wb = $excel.ActiveWorkbook

$updates = @{
    "F2"  = 8347
    "F3"  = 7799
    "F4"  = 126
    "F9"  = 119
    "F10" = 165
    "F11" = 230
    "F14" = 1343
    "F17" = 13
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
